$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.Value = "'2014-04-22"
}
